$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 2.276052666666667
$ws.Cells.Item(2, 8).Value2 = 6.828158
$ws.Cells.Item(2, 9).Value2 = 0.005247614157263819
$ws.Cells.Item(2, 10).Value2 = 0.005247614157263819
$ws.Cells.Item(2, 13).Value2 = 4.277890333333334
$ws.Cells.Item(2, 14).Value2 = 12.833671
$ws.Cells.Item(2, 15).Value2 = 0.04123357425337639
$ws.Cells.Item(2, 16).Value2 = 0.04123357425337638
$ws.Cells.Item(2, 17).Value2 = 9.736703700890892
$ws.Cells.Item(2, 18).Value2 = 87.63033330801802
$ws.Cells.Item(2, 19).Value2 = 0.0002163778880066069
$ws.Cells.Item(2, 20).Value2 = 0.0002163778880066068
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 2.276052666666667
$ws.Cells.Item(3, 8).Value2 = 6.828158
$ws.Cells.Item(3, 9).Value2 = 0.005247614157263819
$ws.Cells.Item(3, 10).Value2 = 0.005247614157263819
$ws.Cells.Item(3, 15).Value2 = 0.4451428460610328
$ws.Cells.Item(3, 16).Value2 = 0.4451428460610327
$ws.Cells.Item(3, 17).Value2 = 105.1139532564936
$ws.Cells.Item(3, 18).Value2 = 946.0255793084419
$ws.Cells.Item(3, 19).Value2 = 0.002335937900994584
$ws.Cells.Item(3, 20).Value2 = 0.002335937900994584
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 2.276052666666667
$ws.Cells.Item(4, 8).Value2 = 6.828158
$ws.Cells.Item(4, 9).Value2 = 0.005247614157263819
$ws.Cells.Item(4, 10).Value2 = 0.005247614157263819
$ws.Cells.Item(4, 13).Value2 = 8.558147333333332
$ws.Cells.Item(4, 14).Value2 = 25.674442
$ws.Cells.Item(4, 15).Value2 = 0.08248996024761777
$ws.Cells.Item(4, 16).Value2 = 0.08248996024761777
$ws.Cells.Item(4, 17).Value2 = 19.47879405975955
$ws.Cells.Item(4, 18).Value2 = 175.309146537836
$ws.Cells.Item(4, 19).Value2 = 0.0004328754832275286
$ws.Cells.Item(4, 20).Value2 = 0.0004328754832275287
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 2.276052666666667
$ws.Cells.Item(5, 8).Value2 = 6.828158
$ws.Cells.Item(5, 9).Value2 = 0.005247614157263819
$ws.Cells.Item(5, 10).Value2 = 0.005247614157263819
$ws.Cells.Item(5, 13).Value2 = 44.72914066666667
$ws.Cells.Item(5, 14).Value2 = 134.187422
$ws.Cells.Item(5, 15).Value2 = 0.4311336194379731
$ws.Cells.Item(5, 16).Value2 = 0.431133619437973
$ws.Cells.Item(5, 17).Value2 = 101.8058798920751
$ws.Cells.Item(5, 18).Value2 = 916.252919028676
$ws.Cells.Item(5, 19).Value2 = 0.002262422885035099
$ws.Cells.Item(5, 20).Value2 = 0.002262422885035099
$ws.Cells.Item(6, 9).Value2 = 0.1062533062835484
$ws.Cells.Item(6, 10).Value2 = 0.1062533062835484
$ws.Cells.Item(6, 13).Value2 = 4.277890333333334
$ws.Cells.Item(6, 14).Value2 = 12.833671
$ws.Cells.Item(6, 15).Value2 = 0.04123357425337639
$ws.Cells.Item(6, 16).Value2 = 0.04123357425337638
$ws.Cells.Item(6, 17).Value2 = 197.1480618655758
$ws.Cells.Item(6, 18).Value2 = 1774.332556790182
$ws.Cells.Item(6, 19).Value2 = 0.004381203594309437
$ws.Cells.Item(6, 20).Value2 = 0.004381203594309436
$ws.Cells.Item(7, 9).Value2 = 0.1062533062835484
$ws.Cells.Item(7, 10).Value2 = 0.1062533062835484
$ws.Cells.Item(7, 15).Value2 = 0.4451428460610328
$ws.Cells.Item(7, 16).Value2 = 0.4451428460610327
$ws.Cells.Item(7, 19).Value2 = 0.04729789916245335
$ws.Cells.Item(7, 20).Value2 = 0.04729789916245335
$ws.Cells.Item(8, 9).Value2 = 0.1062533062835484
$ws.Cells.Item(8, 10).Value2 = 0.1062533062835484
$ws.Cells.Item(8, 13).Value2 = 8.558147333333332
$ws.Cells.Item(8, 14).Value2 = 25.674442
$ws.Cells.Item(8, 15).Value2 = 0.08248996024761777
$ws.Cells.Item(8, 16).Value2 = 0.08248996024761777
$ws.Cells.Item(8, 17).Value2 = 394.4051923865071
$ws.Cells.Item(8, 18).Value2 = 3549.646731478564
$ws.Cells.Item(8, 19).Value2 = 0.008764831011507861
$ws.Cells.Item(8, 20).Value2 = 0.008764831011507863
$ws.Cells.Item(9, 9).Value2 = 0.1062533062835484
$ws.Cells.Item(9, 10).Value2 = 0.1062533062835484
$ws.Cells.Item(9, 13).Value2 = 44.72914066666667
$ws.Cells.Item(9, 14).Value2 = 134.187422
$ws.Cells.Item(9, 15).Value2 = 0.4311336194379731
$ws.Cells.Item(9, 16).Value2 = 0.431133619437973
$ws.Cells.Item(9, 17).Value2 = 2061.357983544859
$ws.Cells.Item(9, 18).Value2 = 18552.22185190372
$ws.Cells.Item(9, 19).Value2 = 0.04580937251527775
$ws.Cells.Item(9, 20).Value2 = 0.04580937251527775
$ws.Cells.Item(10, 7).Value2 = 41.187613
$ws.Cells.Item(10, 8).Value2 = 123.562839
$ws.Cells.Item(10, 9).Value2 = 0.09496120377532416
$ws.Cells.Item(10, 10).Value2 = 0.09496120377532417
$ws.Cells.Item(10, 13).Value2 = 4.277890333333334
$ws.Cells.Item(10, 14).Value2 = 12.833671
$ws.Cells.Item(10, 15).Value2 = 0.04123357425337639
$ws.Cells.Item(10, 16).Value2 = 0.04123357425337638
$ws.Cells.Item(10, 17).Value2 = 176.1960915057744
$ws.Cells.Item(10, 18).Value2 = 1585.764823551969
$ws.Cells.Item(10, 19).Value2 = 0.003915589847059835
$ws.Cells.Item(10, 20).Value2 = 0.003915589847059835
$ws.Cells.Item(11, 7).Value2 = 41.187613
$ws.Cells.Item(11, 8).Value2 = 123.562839
$ws.Cells.Item(11, 9).Value2 = 0.09496120377532416
$ws.Cells.Item(11, 10).Value2 = 0.09496120377532417
$ws.Cells.Item(11, 15).Value2 = 0.4451428460610328
$ws.Cells.Item(11, 16).Value2 = 0.4451428460610327
$ws.Cells.Item(11, 17).Value2 = 1902.149669484162
$ws.Cells.Item(11, 18).Value2 = 17119.34702535746
$ws.Cells.Item(11, 19).Value2 = 0.04227130051392949
$ws.Cells.Item(11, 20).Value2 = 0.04227130051392949
$ws.Cells.Item(12, 7).Value2 = 41.187613
$ws.Cells.Item(12, 8).Value2 = 123.562839
$ws.Cells.Item(12, 9).Value2 = 0.09496120377532416
$ws.Cells.Item(12, 10).Value2 = 0.09496120377532417
$ws.Cells.Item(12, 13).Value2 = 8.558147333333332
$ws.Cells.Item(12, 14).Value2 = 25.674442
$ws.Cells.Item(12, 15).Value2 = 0.08248996024761777
$ws.Cells.Item(12, 16).Value2 = 0.08248996024761777
$ws.Cells.Item(12, 17).Value2 = 352.4896603623153
$ws.Cells.Item(12, 18).Value2 = 3172.406943260838
$ws.Cells.Item(12, 19).Value2 = 0.00783334592449242
$ws.Cells.Item(12, 20).Value2 = 0.007833345924492422
$ws.Cells.Item(13, 7).Value2 = 41.187613
$ws.Cells.Item(13, 8).Value2 = 123.562839
$ws.Cells.Item(13, 9).Value2 = 0.09496120377532416
$ws.Cells.Item(13, 10).Value2 = 0.09496120377532417
$ws.Cells.Item(13, 13).Value2 = 44.72914066666667
$ws.Cells.Item(13, 14).Value2 = 134.187422
$ws.Cells.Item(13, 15).Value2 = 0.4311336194379731
$ws.Cells.Item(13, 16).Value2 = 0.431133619437973
$ws.Cells.Item(13, 17).Value2 = 1842.286535601228
$ws.Cells.Item(13, 18).Value2 = 16580.57882041106
$ws.Cells.Item(13, 19).Value2 = 0.04094096748984242
$ws.Cells.Item(13, 20).Value2 = 0.04094096748984242
$ws.Cells.Item(14, 7).Value2 = 344.1819356666667
$ws.Cells.Item(14, 8).Value2 = 1032.545807
$ws.Cells.Item(14, 9).Value2 = 0.7935378757838636
$ws.Cells.Item(14, 10).Value2 = 0.7935378757838637
$ws.Cells.Item(14, 13).Value2 = 4.277890333333334
$ws.Cells.Item(14, 14).Value2 = 12.833671
$ws.Cells.Item(14, 15).Value2 = 0.04123357425337639
$ws.Cells.Item(14, 16).Value2 = 0.04123357425337638
$ws.Cells.Item(14, 17).Value2 = 1472.372575496389
$ws.Cells.Item(14, 18).Value2 = 13251.3531794675
$ws.Cells.Item(14, 19).Value2 = 0.03272040292400051
$ws.Cells.Item(14, 20).Value2 = 0.03272040292400051
$ws.Cells.Item(15, 7).Value2 = 344.1819356666667
$ws.Cells.Item(15, 8).Value2 = 1032.545807
$ws.Cells.Item(15, 9).Value2 = 0.7935378757838636
$ws.Cells.Item(15, 10).Value2 = 0.7935378757838637
$ws.Cells.Item(15, 15).Value2 = 0.4451428460610328
$ws.Cells.Item(15, 16).Value2 = 0.4451428460610327
$ws.Cells.Item(15, 17).Value2 = 15895.2050746609
$ws.Cells.Item(15, 18).Value2 = 143056.8456719481
$ws.Cells.Item(15, 19).Value2 = 0.3532377084836553
$ws.Cells.Item(15, 20).Value2 = 0.3532377084836553
$ws.Cells.Item(16, 7).Value2 = 344.1819356666667
$ws.Cells.Item(16, 8).Value2 = 1032.545807
$ws.Cells.Item(16, 9).Value2 = 0.7935378757838636
$ws.Cells.Item(16, 10).Value2 = 0.7935378757838637
$ws.Cells.Item(16, 13).Value2 = 8.558147333333332
$ws.Cells.Item(16, 14).Value2 = 25.674442
$ws.Cells.Item(16, 15).Value2 = 0.08248996024761777
$ws.Cells.Item(16, 16).Value2 = 0.08248996024761777
$ws.Cells.Item(16, 17).Value2 = 2945.559714907188
$ws.Cells.Item(16, 18).Value2 = 26510.03743416469
$ws.Cells.Item(16, 19).Value2 = 0.06545890782838996
$ws.Cells.Item(16, 20).Value2 = 0.06545890782838996
$ws.Cells.Item(17, 7).Value2 = 344.1819356666667
$ws.Cells.Item(17, 8).Value2 = 1032.545807
$ws.Cells.Item(17, 9).Value2 = 0.7935378757838636
$ws.Cells.Item(17, 10).Value2 = 0.7935378757838637
$ws.Cells.Item(17, 13).Value2 = 44.72914066666667
$ws.Cells.Item(17, 14).Value2 = 134.187422
$ws.Cells.Item(17, 15).Value2 = 0.4311336194379731
$ws.Cells.Item(17, 16).Value2 = 0.431133619437973
$ws.Cells.Item(17, 17).Value2 = 15394.96221535995
$ws.Cells.Item(17, 18).Value2 = 138554.6599382396
$ws.Cells.Item(17, 19).Value2 = 0.3421208565478178
$ws.Cells.Item(17, 20).Value2 = 0.3421208565478178
